# Expand Swedish party abbreviations to their full names everywhere they
# appear in the workbook (single codes like "C -" as well as "+"-joined
# coalition labels like "C -+KD ").
$wb = $excel.ActiveWorkbook

$map = @{
    "C -" = "C - Centre Party (Centerpartiet, C)";
    "KD " = "KD - Christian Democrats (Kristdemokraterna, KD)";
    "L -" = "L - Liberals (Liberalerna, L), known until  as People's Party Liberals (FP, Folkpartiet liberalerna)";
    "M -" = "M - Moderate Party (Moderata samlingspartiet, M)";
    "NYD" = "NYD - New Democracy (Ny Demokrati, NYD)";
    "S -" = "S - Social Democrats (Socialdemokraterna, S)";
    "V -" = "V - Left Party (Vänsterpartiet, V)";
    "MP " = "MP - Green Party (Miljöpartiet de gröna, MP)";
    "SD " = "SD - Sweden Democrats (Sverigedemokraterna, SD)";
}

function Convert-PartyLabel($v) {
    if ($v -eq $null) { return $null }
    if (-not ($v -is [string])) { return $null }
    if ($v.Contains('+')) {
        $parts = $v.Split('+')
        $changed = $false
        $newParts = @()
        foreach ($p in $parts) {
            if ($map.ContainsKey($p)) {
                $newParts += $map[$p]
                $changed = $true
            } else {
                $newParts += $p
            }
        }
        if ($changed) { return [string]::Join('+', $newParts) }
        return $null
    } elseif ($map.ContainsKey($v)) {
        return $map[$v]
    }
    return $null
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column
    $vals = $used.Value2

    if ($rows -eq 1 -and $cols -eq 1) {
        $newVal = Convert-PartyLabel $vals
        if ($newVal -ne $null) {
            $ws.Cells.Item($startRow, $startCol).Value2 = $newVal
        }
        continue
    }

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $v = $vals[$r, $c]
            $newVal = Convert-PartyLabel $v
            if ($newVal -ne $null) {
                $ws.Cells.Item($startRow + $r - 1, $startCol + $c - 1).Value2 = $newVal
            }
        }
    }
}
